# Adds a "ThirdPass" worksheet after the existing "SecondPass" sheet and
# populates it with a flattened / enriched view of the FirstPass/SecondPass
# product data (identification columns + numeric measures).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ThirdPass"

# Header row
$headers = @("varenavn oplyst","vare","ingrediens","kategori","Vare Nr.","Masse per styk","antal","Samlet masse","tabel")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows: varenavn oplyst, vare, ingrediens, kategori, Vare Nr., Masse per styk, antal, Samlet masse, tabel
$rows = @(
    @("Rodfrugtmix - Tern 20x20mm (3kg)", "Frugtmix", "frugt, blandet", "frugt og bær", "3932-3", 3, 1, 3, 1),
    @("Rodfrugtmix - Tern 20x20mm (5kg)", "Frugtmix", "frugt, blandet", "frugt og bær", "3932-5", 5, 5, 25, 1),
    @("Rødløg - Tern 5x5mm (1kg)", "rødløg", "løg", "grøntsager", "4231-1", 1, 3, 3, 1),
    @("Porre - Skiver 2mm (1kg)", "porre", "porre", "grøntsager", "4441-1", 1, 3, 3, 1),
    @("Frugtsalat m/druer , håndskåret i lage (3,2kg)", "DRUER", "vindrue", "frugt og bær", "7505-32", 3.2, 9, 28.8, 1),
    @("Frugtblanding (U/druer), 15x15mm SMÅ TERN håndskåret (1kg)", "DRUER", "vindrue", "frugt og bær", "7508-1", 1, 74, 74, 1),
    @("Kartoffel 15-25mm u/skræl (3kg)", "Kartoffel", "kartoffel", "grøntsager", "8000-3", 3, 13, 39, 1),
    @("Kartoffel 25-40mm u/skræl (3kg)", "Kartoffel", "kartoffel", "grøntsager", "8005-3", 3, 4, 12, 1),
    @("Kartoffelbåde m/skræl (3kg Vakuum) - Økologisk", "Kartoffel", "kartoffel", "grøntsager", "Ø8551-3", 3, 11, 33, 2),
    @("ØKO Kartoffel - forkogt, mos (3kg)", "Kartoffel", "kartoffel", "grøntsager", "Ø8607-3", 3, 29, 87, 2)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $row[$c]
    }
}
